$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14, Q1=15, matching the formatting of the rest of row 1 ---
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# --- Data rows (2-25): swap I<->K and M<->O values, then fill new P & Q columns with 2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value2  = $kVal  # I = old K
    $ws.Cells.Item($r, 11).Value2 = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value2 = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value2 = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value2 = 2      # P
    $ws.Cells.Item($r, 17).Value2 = 2      # Q
}
